$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47, shifting existing rows 47:73 down to 48:74.
$ws.Rows("47:47").Insert()

# Populate the newly inserted row 47 with the new record.
$ws.Range("A47").Value = 6
$ws.Range("B47").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 44663
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = 100114007
$ws.Range("G47").Value = "Jengibre"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 150
$ws.Range("K47").Value = 10000
$ws.Range("L47").Value = 10000
$ws.Range("M47").Value = 10000
$ws.Range("N47").Value = "$/caja 13 kilos"
$ws.Range("O47").Value = "Perú"
$ws.Range("P47").Value = 769
$ws.Range("Q47").Value = 13
$ws.Range("R47").Value = "Hortaliza"
